$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row right after the header data row (row 2), shifting the
# existing "BlockCompanyBrand" row (and everything below it) down by one.
[void]$ws.Range("A3").EntireRow.Insert()

# Populate the new row 3 with the Company Vendor mapping: UnitName,
# TestCases and Description columns, matching the "Data Sheet And Name"
# mapping style used by the surrounding rows.
$ws.Range("A3").Value = "CompanyVendor"
$ws.Range("B3").Value = "CreateCompanyVendor"
$ws.Range("C3").Value = "CreateCompanyVendor"
$ws.Range("D3:I3").Value = "No"

# Restore the view: select C7 and let the top-left cell follow naturally
# (clears the previous topLeftCell/selection scrolled to C18).
[void]$ws.Activate()
[void]$ws.Range("C7").Select()
